$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Remove the old "LOTO (SOPs)" row (row 3); subsequent rows shift up one.
$ws.Rows.Item(3).Delete()

# Fix the SN (column A) numbering for the rows that shifted up.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2

# Update the "period to expire" figures for the new progress date (04-Nov-2025).
$ws.Cells.Item(3, 8).Value = -229
$ws.Cells.Item(4, 8).Value = 106

# "LAST UPDATE" column (I) holds date-like text, not a real Excel date, in the
# source file. Force text entry (bypassing Excel's date auto-detection) while
# restoring the row's original cell formatting (fill/border/alignment) via a
# formats-only paste from a cell on the same row that still has that style.
$i3 = $ws.Cells.Item(3, 9)
$i3.NumberFormat = "@"
$i3.Value = "04-Nov-2025"
$ws.Cells.Item(3, 8).Copy()
$i3.PasteSpecial(-4122)

$i4 = $ws.Cells.Item(4, 9)
$i4.NumberFormat = "@"
$i4.Value = "04-Nov-2025"
$ws.Cells.Item(4, 8).Copy()
$i4.PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Column C width changed from 19 to 10 (ColumnWidth uses a different unit than
# the stored <col width>, which runs ~0.8333 higher once round-tripped).
$ws.Columns.Item(3).ColumnWidth = 10 - 5/6
